$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.238.09'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '1.848.13'
$ws.Range("E3").Value = '  +0.86%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9995'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '241.02'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.6729'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.54%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.07429'
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.2963'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.58%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '22.92'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.80%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07724'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.87%  '
$ws.Range("D12").Value = '1.831.45'
$ws.Range("E12").Value = '  -0.17%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.018'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.83%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.6787'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.40%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '86.27'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.82%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '6.160'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").Value = '29.175.14'
$ws.Range("E17").Value = '  +0.04%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000008312'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.66%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '228.82'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.23%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '12.56'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("E21").Value = '  +0.03%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '7.207'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.85%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +0.66%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '8.695'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.53%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.1407'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -3.24%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '18.03'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.17%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.508'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.11%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '4.182'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.72%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '4.081'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -1.42%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.193'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.58%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.05324'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +3.33%  '
$ws.Range("E33").Value = '  +2.96%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.7580'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("E35").Value = '  +0.81%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.686'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("D37").Value = '1.333.58'
$ws.Range("E37").Value = '  +1.94%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01804'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.61%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.735'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.40%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.9258'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.36%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '5.956'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +2.75%  '
$ws.Range("E42").Value = '  +0.14%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '103.55'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.70%  '
$ws.Range("B44").Value = 'XinFinNetwork'
$ws.Range("C44").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.07932'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +9.79%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.00000000124'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.971.91'
$ws.Range("E46").Value = '  -0.56%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.5162'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.73%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.771'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.05%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '63.99'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.30%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '9.224'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -3.24%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.05942'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.38%  '
